$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy formatting (bold / border / center-top alignment) from the
# previous year's row-header cell so the new row header matches.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = "2021年"

$values = @(2678, 22043, 5351, 42851, 625, 731, 4507, 3885, 44601, 453, 2662, 139074, 703, 1152, 2298, 4156, 378)

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = 2 + $i  # Column B is 2
    $ws.Cells.Item(10, $col).Value = $values[$i]
}
